# Insert two new price records (rows 246 and 247) into the daily-price
# dataset for Hortaliza / Terminal Hortofrutícola Agro Chillán - Pimiento.
# All existing rows from row 246 onward shift down by two rows
# (old 246 -> new 248, ..., old 313 -> new 315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 246:313 down by inserting two blank rows at 246.
$ws.Rows("246:247").Insert()

# --- New row 246: Zafiro rojo, Región de Arica y Parinacota ---
$ws.Range("A246").Value = 7
$ws.Range("B246").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C246").Value = "Ñuble"
$ws.Range("D246").Value = 44798
$ws.Range("E246").Value = 16
$ws.Range("F246").Value = 100112002
$ws.Range("G246").Value = "Pimiento"
$ws.Range("H246").Value = "Zafiro rojo"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 120
$ws.Range("K246").Value = 26000
$ws.Range("L246").Value = 27000
$ws.Range("M246").Value = 26500
$ws.Range("N246").Value = "`$/caja 15 kilos"
$ws.Range("O246").Value = "Región de Arica y Parinacota"
$ws.Range("P246").Value = 1767
$ws.Range("Q246").Value = 15
$ws.Range("R246").Value = "Hortaliza"

# --- New row 247: Zafiro verde, Región de Arica y Parinacota ---
$ws.Range("A247").Value = 7
$ws.Range("B247").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C247").Value = "Ñuble"
$ws.Range("D247").Value = 44798
$ws.Range("E247").Value = 16
$ws.Range("F247").Value = 100112002
$ws.Range("G247").Value = "Pimiento"
$ws.Range("H247").Value = "Zafiro verde"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 120
$ws.Range("K247").Value = 25000
$ws.Range("L247").Value = 26000
$ws.Range("M247").Value = 25500
$ws.Range("N247").Value = "`$/caja 15 kilos"
$ws.Range("O247").Value = "Región de Arica y Parinacota"
$ws.Range("P247").Value = 1700
$ws.Range("Q247").Value = 15
$ws.Range("R247").Value = "Hortaliza"
